$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(53, 8).Value = 445.55554
$ws.Cells.Item(53, 9).Value = 319.2
$ws.Cells.Item(53, 10).Value = 603.5
$ws.Cells.Item(53, 11).Value = 319.2
$ws.Cells.Item(53, 12).Value = 603.5
$ws.Cells.Item(53, 13).Value = 317.8
$ws.Cells.Item(53, 14).Value = -1877.5

$ws.Cells.Item(64, 8).Value = 8870
$ws.Cells.Item(64, 9).Value = 4187.5
$ws.Cells.Item(64, 10).Value = 11991.667
$ws.Cells.Item(64, 11).Value = 4187.5
$ws.Cells.Item(64, 12).Value = 11991.667
$ws.Cells.Item(64, 13).Value = -3939.5
$ws.Cells.Item(64, 14).Value = -12487.667

$ws.Cells.Item(67, 8).Value = 8870
$ws.Cells.Item(67, 9).Value = 4187.5
$ws.Cells.Item(67, 10).Value = 11991.667
$ws.Cells.Item(67, 11).Value = 4187.5
$ws.Cells.Item(67, 12).Value = 11991.667
$ws.Cells.Item(67, 13).Value = -3329.5
$ws.Cells.Item(67, 14).Value = -13707.667

$ws.Cells.Item(137, 8).Value = 1478.9722
$ws.Cells.Item(137, 9).Value = 1203.8572
$ws.Cells.Item(137, 11).Value = 3611.5716
$ws.Cells.Item(137, 13).Value = -1061.5716

$ws.Cells.Item(138, 8).Value = 6175342.5
$ws.Cells.Item(138, 9).Value = 1161.7241
$ws.Cells.Item(138, 10).Value = 9618636
$ws.Cells.Item(138, 11).Value = 3485.1723
$ws.Cells.Item(138, 12).Value = 28855908
$ws.Cells.Item(138, 13).Value = 1654.8277
$ws.Cells.Item(138, 14).Value = -28866188

$ws.Cells.Item(139, 8).Value = 76299.8
$ws.Cells.Item(139, 10).Value = 76624.75
$ws.Cells.Item(139, 12).Value = 76624.75
$ws.Cells.Item(139, 14).Value = -86904.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(88, 8).Value = 2233.6667
$ws.Cells.Item(88, 9).Value = 2553
$ws.Cells.Item(88, 10).Value = 2074
$ws.Cells.Item(88, 11).Value = 2553
$ws.Cells.Item(88, 12).Value = 2074
$ws.Cells.Item(88, 13).Value = -2147
$ws.Cells.Item(88, 14).Value = -2886

$ws.Cells.Item(91, 8).Value = 2233.6667
$ws.Cells.Item(91, 9).Value = 2553
$ws.Cells.Item(91, 10).Value = 2074
$ws.Cells.Item(91, 11).Value = 2553
$ws.Cells.Item(91, 12).Value = 2074
$ws.Cells.Item(91, 13).Value = -1149
$ws.Cells.Item(91, 14).Value = -4882

$ws.Cells.Item(97, 8).Value = 2385.4375
$ws.Cells.Item(97, 10).Value = 3824.75
$ws.Cells.Item(97, 12).Value = 3824.75
$ws.Cells.Item(97, 14).Value = -4816.75

$ws.Cells.Item(105, 8).Value = 65821.5
$ws.Cells.Item(105, 10).Value = 65821.5
$ws.Cells.Item(105, 12).Value = 65821.5
$ws.Cells.Item(105, 14).Value = -72809.5

$ws.Cells.Item(106, 8).Value = 30000
$ws.Cells.Item(106, 10).Value = 30000
$ws.Cells.Item(106, 12).Value = 30000
$ws.Cells.Item(106, 14).Value = -32524

$ws.Cells.Item(122, 8).Value = 1101.2106
$ws.Cells.Item(122, 9).Value = 838.1875
$ws.Cells.Item(122, 10).Value = 2504
$ws.Cells.Item(122, 11).Value = 2514.5625
$ws.Cells.Item(122, 12).Value = 7512
$ws.Cells.Item(122, 13).Value = -64.5625
$ws.Cells.Item(122, 14).Value = -12412

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 2408.2
$ws.Cells.Item(86, 9).Value = 2351.25
$ws.Cells.Item(86, 10).Value = 2636
$ws.Cells.Item(86, 11).Value = 2351.25
$ws.Cells.Item(86, 12).Value = 2636
$ws.Cells.Item(86, 13).Value = -1228.25
$ws.Cells.Item(86, 14).Value = -4882

$ws.Cells.Item(89, 8).Value = 2408.2
$ws.Cells.Item(89, 9).Value = 2351.25
$ws.Cells.Item(89, 10).Value = 2636
$ws.Cells.Item(89, 11).Value = 11756.25
$ws.Cells.Item(89, 12).Value = 13180
$ws.Cells.Item(89, 13).Value = -6140.25
$ws.Cells.Item(89, 14).Value = -24412

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(80, 8).Value = 40128
$ws.Cells.Item(80, 10).Value = 40128
$ws.Cells.Item(80, 12).Value = 40128
$ws.Cells.Item(80, 14).Value = -42374

$ws.Cells.Item(83, 8).Value = 40128
$ws.Cells.Item(83, 10).Value = 40128
$ws.Cells.Item(83, 12).Value = 120384
$ws.Cells.Item(83, 14).Value = -131616

$ws.Cells.Item(99, 8).Value = 3827.5
$ws.Cells.Item(99, 10).Value = 4254
$ws.Cells.Item(99, 12).Value = 4254
$ws.Cells.Item(99, 14).Value = -7250

$ws.Cells.Item(126, 8).Value = 3827.5
$ws.Cells.Item(126, 10).Value = 4254
$ws.Cells.Item(126, 12).Value = 12762
$ws.Cells.Item(126, 14).Value = -17702

$ws.Cells.Item(131, 8).Value = 31824.75
$ws.Cells.Item(131, 10).Value = 31824.75
$ws.Cells.Item(131, 12).Value = 31824.75
$ws.Cells.Item(131, 14).Value = -41904.75

$ws.Cells.Item(132, 8).Value = 2544.963
$ws.Cells.Item(132, 9).Value = 2459.3726
$ws.Cells.Item(132, 11).Value = 7378.1178
$ws.Cells.Item(132, 13).Value = -4848.1178

$ws.Cells.Item(141, 8).Value = 372084.16
$ws.Cells.Item(141, 10).Value = 372084.16
$ws.Cells.Item(141, 12).Value = 372084.16
$ws.Cells.Item(141, 14).Value = -382444.16

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(32, 8).Value = 4391.5
$ws.Cells.Item(32, 9).Value = 4500
$ws.Cells.Item(32, 10).Value = 4283
$ws.Cells.Item(32, 11).Value = 13500
$ws.Cells.Item(32, 12).Value = 12849
$ws.Cells.Item(32, 13).Value = -13217
$ws.Cells.Item(32, 14).Value = -13415

$ws.Cells.Item(46, 8).Value = 501850
$ws.Cells.Item(46, 9).Value = 501850
$ws.Cells.Item(46, 11).Value = 1505550
$ws.Cells.Item(46, 13).Value = -1505459

$ws.Cells.Item(138, 8).Value = 27788152
$ws.Cells.Item(138, 9).Value = 38471324
$ws.Cells.Item(138, 10).Value = 11900
$ws.Cells.Item(138, 11).Value = 115413972
$ws.Cells.Item(138, 12).Value = 35700
$ws.Cells.Item(138, 13).Value = -115408832
$ws.Cells.Item(138, 14).Value = -45980

$ws.Cells.Item(139, 8).Value = 1673.2727
$ws.Cells.Item(139, 9).Value = 1378.4445
$ws.Cells.Item(139, 11).Value = 4135.333500000001
$ws.Cells.Item(139, 13).Value = 1004.666499999999

$ws.Cells.Item(140, 8).Value = 4071
$ws.Cells.Item(140, 9).Value = 3899.4
$ws.Cells.Item(140, 11).Value = 11698.2
$ws.Cells.Item(140, 13).Value = -6518.200000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(26, 8).Value = 21250
$ws.Cells.Item(26, 10).Value = 21250
$ws.Cells.Item(26, 12).Value = 21250
$ws.Cells.Item(26, 14).Value = -21810

$ws.Cells.Item(50, 8).Value = 21250
$ws.Cells.Item(50, 10).Value = 21250
$ws.Cells.Item(50, 12).Value = 21250
$ws.Cells.Item(50, 14).Value = -22246

$ws.Cells.Item(52, 8).Value = 57649.25
$ws.Cells.Item(52, 10).Value = 57649.25
$ws.Cells.Item(52, 12).Value = 57649.25
$ws.Cells.Item(52, 14).Value = -58167.25

$ws.Cells.Item(58, 8).Value = 31675
$ws.Cells.Item(58, 9).Value = 27100
$ws.Cells.Item(58, 10).Value = 33200
$ws.Cells.Item(58, 11).Value = 27100
$ws.Cells.Item(58, 12).Value = 33200
$ws.Cells.Item(58, 13).Value = -26823
$ws.Cells.Item(58, 14).Value = -33754

$ws.Cells.Item(97, 8).Value = 1087.4445
$ws.Cells.Item(97, 9).Value = 1173.375
$ws.Cells.Item(97, 11).Value = 1173.375
$ws.Cells.Item(97, 13).Value = -677.375

$ws.Cells.Item(122, 8).Value = 2600.5862
$ws.Cells.Item(122, 9).Value = 2169.9092
$ws.Cells.Item(122, 10).Value = 3954.1428
$ws.Cells.Item(122, 11).Value = 6509.7276
$ws.Cells.Item(122, 12).Value = 11862.4284
$ws.Cells.Item(122, 13).Value = -4059.7276
$ws.Cells.Item(122, 14).Value = -16762.4284

$ws.Cells.Item(132, 8).Value = 2946.25
$ws.Cells.Item(132, 9).Value = 2331.4211
$ws.Cells.Item(132, 11).Value = 6994.263300000001
$ws.Cells.Item(132, 13).Value = -4464.263300000001

$ws.Cells.Item(139, 8).Value = 125000
$ws.Cells.Item(139, 10).Value = 125000
$ws.Cells.Item(139, 12).Value = 125000
$ws.Cells.Item(139, 14).Value = -135280

$ws.Cells.Item(141, 8).Value = 111103.5
$ws.Cells.Item(141, 10).Value = 111103.5
$ws.Cells.Item(141, 12).Value = 111103.5
$ws.Cells.Item(141, 14).Value = -121463.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 5587.4375
$ws.Cells.Item(7, 9).Value = 5463.909
$ws.Cells.Item(7, 10).Value = 5859.2
$ws.Cells.Item(7, 11).Value = 5463.909
$ws.Cells.Item(7, 12).Value = 5859.2
$ws.Cells.Item(7, 13).Value = -5351.909
$ws.Cells.Item(7, 14).Value = -6083.2

$ws.Cells.Item(80, 8).Value = 39862.5
$ws.Cells.Item(80, 10).Value = 39862.5
$ws.Cells.Item(80, 12).Value = 39862.5
$ws.Cells.Item(80, 14).Value = -42108.5

$ws.Cells.Item(83, 8).Value = 39862.5
$ws.Cells.Item(83, 10).Value = 39862.5
$ws.Cells.Item(83, 12).Value = 119587.5
$ws.Cells.Item(83, 14).Value = -130819.5

$ws.Cells.Item(122, 8).Value = 255058.05
$ws.Cells.Item(122, 9).Value = 316988.22
$ws.Cells.Item(122, 10).Value = 7337.375
$ws.Cells.Item(122, 11).Value = 950964.6599999999
$ws.Cells.Item(122, 12).Value = 22012.125
$ws.Cells.Item(122, 13).Value = -948514.6599999999
$ws.Cells.Item(122, 14).Value = -26912.125

$ws.Cells.Item(126, 8).Value = 5587.4375
$ws.Cells.Item(126, 9).Value = 5463.909
$ws.Cells.Item(126, 10).Value = 5859.2
$ws.Cells.Item(126, 11).Value = 16391.727
$ws.Cells.Item(126, 12).Value = 17577.6
$ws.Cells.Item(126, 13).Value = -13921.727
$ws.Cells.Item(126, 14).Value = -22517.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 7248.875
$ws.Cells.Item(62, 9).Value = 4998.3335
$ws.Cells.Item(62, 10).Value = 7768.231
$ws.Cells.Item(62, 11).Value = 4998.3335
$ws.Cells.Item(62, 12).Value = 7768.231
$ws.Cells.Item(62, 13).Value = -4374.3335
$ws.Cells.Item(62, 14).Value = -9016.231

$ws.Cells.Item(65, 8).Value = 7248.875
$ws.Cells.Item(65, 9).Value = 4998.3335
$ws.Cells.Item(65, 10).Value = 7768.231
$ws.Cells.Item(65, 11).Value = 24991.6675
$ws.Cells.Item(65, 12).Value = 38841.155
$ws.Cells.Item(65, 13).Value = -21871.6675
$ws.Cells.Item(65, 14).Value = -45081.155

$ws.Cells.Item(113, 8).Value = 1696.1305
$ws.Cells.Item(113, 9).Value = 1647.7858
$ws.Cells.Item(113, 10).Value = 1771.3334
$ws.Cells.Item(113, 11).Value = 4943.357400000001
$ws.Cells.Item(113, 12).Value = 5314.0002
$ws.Cells.Item(113, 13).Value = -2773.357400000001
$ws.Cells.Item(113, 14).Value = -9654.0002

$ws.Cells.Item(135, 8).Value = 74353.5
$ws.Cells.Item(135, 10).Value = 74353.5
$ws.Cells.Item(135, 12).Value = 74353.5
$ws.Cells.Item(135, 14).Value = -84493.5
